# Update res_bus vm_pu results for the 380 kV case (rows 2-25, columns B-F and I-N).
# Column B (bus 0, slack) is lowered from 1.05 p.u. to 1.02 p.u.; all other
# bus voltage magnitudes are updated to the corresponding recalculated values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.031084389381527
$ws.Cells.Item(2, 4).Value = 1.034571803276337
$ws.Cells.Item(2, 5).Value = 1.044694289852873
$ws.Cells.Item(2, 6).Value = 1.052361646099872
$ws.Cells.Item(2, 9).Value = 1.033315055699082
$ws.Cells.Item(2, 10).Value = 1.03622206645256
$ws.Cells.Item(2, 11).Value = 1.037370893434816
$ws.Cells.Item(2, 12).Value = 1.047464619512056
$ws.Cells.Item(2, 13).Value = 1.055110590174108
$ws.Cells.Item(2, 14).Value = 1.016037272982404

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.032031339333612
$ws.Cells.Item(3, 4).Value = 1.035265900666878
$ws.Cells.Item(3, 5).Value = 1.045627420191635
$ws.Cells.Item(3, 6).Value = 1.053368841898299
$ws.Cells.Item(3, 9).Value = 1.03346621364555
$ws.Cells.Item(3, 10).Value = 1.036810755430282
$ws.Cells.Item(3, 11).Value = 1.037874512993238
$ws.Cells.Item(3, 12).Value = 1.048208721972835
$ws.Cells.Item(3, 13).Value = 1.055930117490809
$ws.Cells.Item(3, 14).Value = 1.016235114999597

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.032644579077515
$ws.Cells.Item(4, 4).Value = 1.035715206614894
$ws.Cells.Item(4, 5).Value = 1.046232075728952
$ws.Cells.Item(4, 6).Value = 1.054021480037969
$ws.Cells.Item(4, 9).Value = 1.033562777048991
$ws.Cells.Item(4, 10).Value = 1.037191568911506
$ws.Cells.Item(4, 11).Value = 1.038199891739502
$ws.Cells.Item(4, 12).Value = 1.048690437003006
$ws.Cells.Item(4, 13).Value = 1.056460708583385
$ws.Cells.Item(4, 14).Value = 1.016363022542629

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.032902503819663
$ws.Cells.Item(5, 4).Value = 1.035904136199157
$ws.Cells.Item(5, 5).Value = 1.046486477049537
$ws.Cells.Item(5, 6).Value = 1.054296066977003
$ws.Cells.Item(5, 9).Value = 1.033603073512884
$ws.Cells.Item(5, 10).Value = 1.037351636228308
$ws.Cells.Item(5, 11).Value = 1.038336561044863
$ws.Cells.Item(5, 12).Value = 1.048893004252951
$ws.Cells.Item(5, 13).Value = 1.056683840304656
$ws.Cells.Item(5, 14).Value = 1.016416768299908

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.032945817456509
$ws.Cells.Item(6, 4).Value = 1.035935860693274
$ws.Cells.Item(6, 5).Value = 1.046529204085442
$ws.Cells.Item(6, 6).Value = 1.054342184072811
$ws.Cells.Item(6, 9).Value = 1.03360982193012
$ws.Cells.Item(6, 10).Value = 1.037378510645054
$ws.Cells.Item(6, 11).Value = 1.038359501361739
$ws.Cells.Item(6, 12).Value = 1.048927019336508
$ws.Cells.Item(6, 13).Value = 1.056721309247235
$ws.Cells.Item(6, 14).Value = 1.016425790880104

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.032648025015039
$ws.Cells.Item(7, 4).Value = 1.035717730940444
$ws.Cells.Item(7, 5).Value = 1.046235474251285
$ws.Cells.Item(7, 6).Value = 1.054025148227543
$ws.Cells.Item(7, 9).Value = 1.033563316667006
$ws.Cells.Item(7, 10).Value = 1.03719370784475
$ws.Cells.Item(7, 11).Value = 1.038201718393891
$ws.Cells.Item(7, 12).Value = 1.048693143504653
$ws.Cells.Item(7, 13).Value = 1.056463689801875
$ws.Cells.Item(7, 14).Value = 1.016363740800944

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.031404311521107
$ws.Cells.Item(8, 4).Value = 1.034806338948401
$ws.Cells.Item(8, 5).Value = 1.04500946726468
$ws.Cells.Item(8, 6).Value = 1.052701842425667
$ws.Cells.Item(8, 9).Value = 1.033366397892203
$ws.Cells.Item(8, 10).Value = 1.036421038398252
$ws.Cells.Item(8, 11).Value = 1.037541196290638
$ws.Cells.Item(8, 12).Value = 1.047716043966893
$ws.Cells.Item(8, 13).Value = 1.055387490049074
$ws.Cells.Item(8, 14).Value = 1.01610415706509

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.02921659709744
$ws.Cells.Item(9, 4).Value = 1.03320177774685
$ws.Cells.Item(9, 5).Value = 1.042855710389526
$ws.Cells.Item(9, 6).Value = 1.050377070773014
$ws.Cells.Item(9, 9).Value = 1.033009882795446
$ws.Cells.Item(9, 10).Value = 1.035058711278959
$ws.Cells.Item(9, 11).Value = 1.036373510786681
$ws.Cells.Item(9, 12).Value = 1.045996084958714
$ws.Cells.Item(9, 13).Value = 1.053493451096408
$ws.Cells.Item(9, 14).Value = 1.015645915080208

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.027760773457715
$ws.Cells.Item(10, 4).Value = 1.032133111539323
$ws.Cells.Item(10, 5).Value = 1.041424398322164
$ws.Cells.Item(10, 6).Value = 1.048832042412607
$ws.Cells.Item(10, 9).Value = 1.032765836730465
$ws.Cells.Item(10, 10).Value = 1.034150018242054
$ws.Cells.Item(10, 11).Value = 1.035592583182028
$ws.Cells.Item(10, 12).Value = 1.044850727457591
$ws.Cells.Item(10, 13).Value = 1.052232411508335
$ws.Cells.Item(10, 14).Value = 1.015339889097257

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.027131025493093
$ws.Cells.Item(11, 4).Value = 1.031670631184579
$ws.Cells.Item(11, 5).Value = 1.040805712778696
$ws.Cells.Item(11, 6).Value = 1.048164185817913
$ws.Cells.Item(11, 9).Value = 1.032658657740583
$ws.Cells.Item(11, 10).Value = 1.033756442344173
$ws.Cells.Item(11, 11).Value = 1.035253858983977
$ws.Cells.Item(11, 12).Value = 1.044355091709872
$ws.Cells.Item(11, 13).Value = 1.051686772837523
$ws.Cells.Item(11, 14).Value = 1.015207254904583

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.026897204998263
$ws.Cells.Item(12, 4).Value = 1.031498885668204
$ws.Cells.Item(12, 5).Value = 1.040576069150937
$ws.Cells.Item(12, 6).Value = 1.04791628818866
$ws.Cells.Item(12, 9).Value = 1.032618621001624
$ws.Cells.Item(12, 10).Value = 1.033610235675869
$ws.Cells.Item(12, 11).Value = 1.035127955928886
$ws.Cells.Item(12, 12).Value = 1.044171038324247
$ws.Cells.Item(12, 13).Value = 1.051484159437166
$ws.Cells.Item(12, 14).Value = 1.015157970540023

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.026947355933865
$ws.Cells.Item(13, 4).Value = 1.031535723820203
$ws.Cells.Item(13, 5).Value = 1.040625321060734
$ws.Cells.Item(13, 6).Value = 1.047969455177351
$ws.Cells.Item(13, 9).Value = 1.032627219223194
$ws.Cells.Item(13, 10).Value = 1.033641598170308
$ws.Cells.Item(13, 11).Value = 1.035154966429401
$ws.Cells.Item(13, 12).Value = 1.044210516221417
$ws.Cells.Item(13, 13).Value = 1.051527617917981
$ws.Cells.Item(13, 14).Value = 1.015168543017898

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.027111695856304
$ws.Cells.Item(14, 4).Value = 1.031656433817173
$ws.Cells.Item(14, 5).Value = 1.040786727017981
$ws.Cells.Item(14, 6).Value = 1.048143690965451
$ws.Cells.Item(14, 9).Value = 1.03265535289245
$ws.Cells.Item(14, 10).Value = 1.033744357157621
$ws.Cells.Item(14, 11).Value = 1.035243453543148
$ws.Cells.Item(14, 12).Value = 1.044339876818812
$ws.Cells.Item(14, 13).Value = 1.051670023498312
$ws.Cells.Item(14, 14).Value = 1.015203181411262

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.02721296384314
$ws.Cells.Item(15, 4).Value = 1.031730812607777
$ws.Cells.Item(15, 5).Value = 1.040886196286103
$ws.Cells.Item(15, 6).Value = 1.04825106648198
$ws.Cells.Item(15, 9).Value = 1.03267265708301
$ws.Cells.Item(15, 10).Value = 1.033807668391342
$ws.Cells.Item(15, 11).Value = 1.035297962032349
$ws.Cells.Item(15, 12).Value = 1.044419586515633
$ws.Cells.Item(15, 13).Value = 1.051757772403144
$ws.Cells.Item(15, 14).Value = 1.015224520877973

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.027802581166361
$ws.Cells.Item(16, 4).Value = 1.032163810412278
$ws.Cells.Item(16, 5).Value = 1.041465481340782
$ws.Cells.Item(16, 6).Value = 1.048876390200191
$ws.Cells.Item(16, 9).Value = 1.032772918170602
$ws.Cells.Item(16, 10).Value = 1.034176136420579
$ws.Cells.Item(16, 11).Value = 1.035615051106765
$ws.Cells.Item(16, 12).Value = 1.044883627818974
$ws.Cells.Item(16, 13).Value = 1.052268632261114
$ws.Cells.Item(16, 14).Value = 1.01534868903166

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.028172602642593
$ws.Cells.Item(17, 4).Value = 1.03243548882958
$ws.Cells.Item(17, 5).Value = 1.041829142064125
$ws.Cells.Item(17, 6).Value = 1.049268948174727
$ws.Cells.Item(17, 9).Value = 1.032835406596262
$ws.Cells.Item(17, 10).Value = 1.034407238944334
$ws.Cells.Item(17, 11).Value = 1.035813798785288
$ws.Cells.Item(17, 12).Value = 1.045174792788958
$ws.Cells.Item(17, 13).Value = 1.052589188916965
$ws.Cells.Item(17, 14).Value = 1.015426543769712

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.028388490957943
$ws.Cells.Item(18, 4).Value = 1.032593979169675
$ws.Cells.Item(18, 5).Value = 1.042041363532879
$ws.Cells.Item(18, 6).Value = 1.049498031566451
$ws.Cells.Item(18, 9).Value = 1.032871709752689
$ws.Cells.Item(18, 10).Value = 1.034542026870204
$ws.Cells.Item(18, 11).Value = 1.035929669170587
$ws.Cells.Item(18, 12).Value = 1.045344654389452
$ws.Cells.Item(18, 13).Value = 1.052776202695727
$ws.Cells.Item(18, 14).Value = 1.015471943237646

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.028462113625177
$ws.Cells.Item(19, 4).Value = 1.032648024465161
$ws.Cells.Item(19, 5).Value = 1.042113743228668
$ws.Cells.Item(19, 6).Value = 1.049576161870249
$ws.Cells.Item(19, 9).Value = 1.032884063533294
$ws.Cells.Item(19, 10).Value = 1.034587984302764
$ws.Cells.Item(19, 11).Value = 1.03596916848894
$ws.Cells.Item(19, 12).Value = 1.045402577864166
$ws.Cells.Item(19, 13).Value = 1.052839976061135
$ws.Cells.Item(19, 14).Value = 1.01548742125674

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.028132896519714
$ws.Cells.Item(20, 4).Value = 1.032406337719675
$ws.Cells.Item(20, 5).Value = 1.041790113911236
$ws.Cells.Item(20, 6).Value = 1.049226818933116
$ws.Cells.Item(20, 9).Value = 1.032828717203581
$ws.Cells.Item(20, 10).Value = 1.034382444886453
$ws.Cells.Item(20, 11).Value = 1.035792480802366
$ws.Cells.Item(20, 12).Value = 1.045143550433102
$ws.Cells.Item(20, 13).Value = 1.052554792231291
$ws.Cells.Item(20, 14).Value = 1.015418191919782

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.027063299190413
$ws.Cells.Item(21, 4).Value = 1.031620886608999
$ws.Cells.Item(21, 5).Value = 1.040739192457369
$ws.Cells.Item(21, 6).Value = 1.048092378054869
$ws.Cells.Item(21, 9).Value = 1.032647074454319
$ws.Cells.Item(21, 10).Value = 1.033714097605904
$ws.Cells.Item(21, 11).Value = 1.03521739865475
$ws.Cells.Item(21, 12).Value = 1.044301782019703
$ws.Cells.Item(21, 13).Value = 1.051628086909402
$ws.Cells.Item(21, 14).Value = 1.015192981765159

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.026391356367238
$ws.Cells.Item(22, 4).Value = 1.031127276028668
$ws.Cells.Item(22, 5).Value = 1.04007938429975
$ws.Cells.Item(22, 6).Value = 1.047380117593882
$ws.Cells.Item(22, 9).Value = 1.032531562916695
$ws.Cells.Item(22, 10).Value = 1.033293794453308
$ws.Cells.Item(22, 11).Value = 1.034855325701058
$ws.Cells.Item(22, 12).Value = 1.04377280584878
$ws.Cells.Item(22, 13).Value = 1.051045784491543
$ws.Cells.Item(22, 14).Value = 1.015051278538647

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.026747513052124
$ws.Cells.Item(23, 4).Value = 1.031388925735975
$ws.Cells.Item(23, 5).Value = 1.040429070850319
$ws.Cells.Item(23, 6).Value = 1.047757604455672
$ws.Cells.Item(23, 9).Value = 1.032592921355573
$ws.Cells.Item(23, 10).Value = 1.033516613042668
$ws.Cells.Item(23, 11).Value = 1.035047314123745
$ws.Cells.Item(23, 12).Value = 1.044053199536341
$ws.Cells.Item(23, 13).Value = 1.051354440120367
$ws.Cells.Item(23, 14).Value = 1.015126407904893

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.028150837818913
$ws.Cells.Item(24, 4).Value = 1.032419509774598
$ws.Cells.Item(24, 5).Value = 1.04180774873134
$ws.Cells.Item(24, 6).Value = 1.049245854981183
$ws.Cells.Item(24, 9).Value = 1.03283174030114
$ws.Cells.Item(24, 10).Value = 1.034393648285747
$ws.Cells.Item(24, 11).Value = 1.035802113653638
$ws.Cells.Item(24, 12).Value = 1.045157667416283
$ws.Cells.Item(24, 13).Value = 1.052570334494061
$ws.Cells.Item(24, 14).Value = 1.015421965798014

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.029781710011505
$ws.Cells.Item(25, 4).Value = 1.033616417689152
$ws.Cells.Item(25, 5).Value = 1.043411715748119
$ws.Cells.Item(25, 6).Value = 1.050977235736923
$ws.Cells.Item(25, 9).Value = 1.033103175387575
$ws.Cells.Item(25, 10).Value = 1.035410993073119
$ws.Cells.Item(25, 11).Value = 1.036675825200942
$ws.Cells.Item(25, 12).Value = 1.046440514158046
$ws.Cells.Item(25, 13).Value = 1.053982818945374
$ws.Cells.Item(25, 14).Value = 1.015764476842926
